# TimeReport_AspenBoler.xlsx -- add new Alpha-Build time-tracking rows
# (rows 37-43) describing ESP32 / PR integration work, matching the
# sharedStrings insertion order so indices line up with the target diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Row 37 ---------------------------------------------------------
$ws.Range("A37").Value = 46030
$ws.Range("C37").Value = 0.59027777777777779
$ws.Range("D37").Value = 0.79166666666666663
$ws.Range("F37").Value = "Experimented with ADC peripheral on DAD and reworked the peripheral  so it only needs to call an initialization function in main, reworked some aspects of main to avoid bad outputs/improve organization,  removed extra plant data structure, created heartbeat module"

# --- Row 38 ---------------------------------------------------------
$ws.Range("A38").Value = 46038
$ws.Range("C38").Value = 0.67222222222222228
$ws.Range("D38").Value = 0.74930555555555556
$ws.Range("F38").Value = "Changed personal directory name and it killed everything, took me an hour but it forced me to reinstall everything to get it working"

# --- Row 39 ---------------------------------------------------------
$ws.Range("A39").Value = 46038
$ws.Range("C39").Value = 0.74930555555555556
$ws.Range("D39").Value = 0.83125000000000004
$ws.Range("F39").Value = "Started modifying main to more closely match a routine based off of received signal from database. Also made some tickets and cleaned up code to be more organized and efficient"

# "Alpha Build" (category) is entered after F37:F39 so the new shared
# string is appended in the same order as the target workbook.
$ws.Range("B37").Value = "Alpha Build"
$ws.Range("B38").Value = "Alpha Build"
$ws.Range("B39").Value = "Alpha Build"

# --- Row 40 ---------------------------------------------------------
$ws.Range("A40").Value = 46043
$ws.Range("B40").Value = "Alpha Build"
$ws.Range("C40").Value = 0.54513888888888884
$ws.Range("D40").Value = 0.58750000000000002
$ws.Range("F40").Value = "Created pin diagram showing which pins were used for what; Integrated fertilizer pump addition to PWM and main and tested it out"

# --- Row 41 ---------------------------------------------------------
$ws.Range("A41").Value = 46045
$ws.Range("B41").Value = "Check-in"
$ws.Range("C41").Value = 0.44444444444444442
$ws.Range("D41").Value = 0.49236111111111114
$ws.Range("F41").Value = "Check-in 6: Met with group & Tyler, discussed what's gotten done and what needs to get done for alpha build, showed bought materials and began observing what modifications would be needed to integrate electronics into plastic build"

# --- Row 42 ---------------------------------------------------------
$ws.Range("A42").Value = 46045
$ws.Range("B42").Value = "Alpha Build"
$ws.Range("C42").Value = 0.72916666666666663
$ws.Range("D42").Value = 0.77152777777777781
$ws.Range("F42").Value = "Realized I should go ahead and have main set up for potential of boolean to change depending on if auto-scheduling/care is on;  Created/modified tickets in project hub."

# --- Row 43 (no Ending Time -> negative calculated duration) --------
$ws.Range("A43").Value = 46048
$ws.Range("B43").Value = "Alpha Build"
$ws.Range("C43").Value = 0.58333333333333337
$ws.Range("F43").Value = "Integrated PR into main; solder and heat shrinked additional wires to extend motor wire length"

# --- Scroll the sheet / move the active selection to match the author's
#     last on-screen position when the file was saved.
$ws.Activate()
$ws.Range("F49").Select()
$excel.ActiveWindow.ScrollRow = 27
$excel.ActiveWindow.ScrollColumn = 3

$wb.Application.Calculate()
